# Generate Report for Handoff
#
# The previous handback run for d6944fb7-...md is dropped (handback finished
# / no longer pending) and the still-open d3a31d31-...md entry moves from
# "Handed back: in sync with en-US" to "Ready for handoff" with refreshed
# handoff timestamps.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview": drop the second data row (d6944fb7...) entirely, and
# refresh the status/date for the remaining row.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Hyperlinks.Delete() on a sub-range removes every hyperlink on the sheet
# in this runtime, so clear them all and re-add only the ones that must
# survive the row-3 deletion.
$wsOverview.Hyperlinks.Delete()

$wsOverview.Rows.Item(3).Delete()

$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"
$wsOverview.Range("D2").Value = "2016-34-11 20:34:29"

$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/49f80eef62279743596beb06eb1b4cd65bf3e36e/e2e/d3a31d31-5aa2-4599-8359-ef1e57e533f9.md", "", "", "d3a31d31-5aa2-4599-8359-ef1e57e533f9.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn": drop the second data row, refresh the status/date on
# the remaining row.
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Hyperlinks.Delete()

$wsZhCn.Rows.Item(3).Delete()

$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("E2").Value = "2016-03-11 20:34:26"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/49f80eef62279743596beb06eb1b4cd65bf3e36e/e2e/d3a31d31-5aa2-4599-8359-ef1e57e533f9.md", "", "", "d3a31d31-5aa2-4599-8359-ef1e57e533f9.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/49f80eef62279743596beb06eb1b4cd65bf3e36e/e2e/d3a31d31-5aa2-4599-8359-ef1e57e533f9.md", "", "", ".md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c06362fa881bd308de4deea933332e8643927c15/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/d3a31d31-5aa2-4599-8359-ef1e57e533f9.4b8fd5364e0172a0c45ee5a376e1f8497edbe73a.zh-cn.xlf", "", "", "d3a31d31-5aa2-4599-8359-ef1e57e533f9.4b8fd5364e0172a0c45ee5a376e1f8497edbe73a.zh-cn.xlf") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/ab779e49edbd22d70415d09974a4aa964bf2125d/e2e/d3a31d31-5aa2-4599-8359-ef1e57e533f9.md", "", "", "d3a31d31-5aa2-4599-8359-ef1e57e533f9.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/73c0dbe6f201bceeb67da08e6ea6e00e75d50143/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/d3a31d31-5aa2-4599-8359-ef1e57e533f9.4b8fd5364e0172a0c45ee5a376e1f8497edbe73a.zh-cn.xlf", "", "", "d3a31d31-5aa2-4599-8359-ef1e57e533f9.4b8fd5364e0172a0c45ee5a376e1f8497edbe73a.zh-cn.xlf") | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de": drop the second data row, refresh the status/date on
# the remaining row.
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Hyperlinks.Delete()

$wsDeDe.Rows.Item(3).Delete()

$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("E2").Value = "2016-03-11 20:34:29"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/49f80eef62279743596beb06eb1b4cd65bf3e36e/e2e/d3a31d31-5aa2-4599-8359-ef1e57e533f9.md", "", "", "d3a31d31-5aa2-4599-8359-ef1e57e533f9.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/49f80eef62279743596beb06eb1b4cd65bf3e36e/e2e/d3a31d31-5aa2-4599-8359-ef1e57e533f9.md", "", "", ".md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d66db1c6ebb25aaa62f08648ab46772a5a8b846a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/d3a31d31-5aa2-4599-8359-ef1e57e533f9.4b8fd5364e0172a0c45ee5a376e1f8497edbe73a.de-de.xlf", "", "", "d3a31d31-5aa2-4599-8359-ef1e57e533f9.4b8fd5364e0172a0c45ee5a376e1f8497edbe73a.de-de.xlf") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/8e983ffe865a07d788c48452891259d6a105a6e3/e2e/d3a31d31-5aa2-4599-8359-ef1e57e533f9.md", "", "", "d3a31d31-5aa2-4599-8359-ef1e57e533f9.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/3c49a4fb6f4e669f928edc86355bbe2e6314ddc3/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/d3a31d31-5aa2-4599-8359-ef1e57e533f9.4b8fd5364e0172a0c45ee5a376e1f8497edbe73a.de-de.xlf", "", "", "d3a31d31-5aa2-4599-8359-ef1e57e533f9.4b8fd5364e0172a0c45ee5a376e1f8497edbe73a.de-de.xlf") | Out-Null
